# Use 3rd quartile instead of mean: update the Num_Inclusions (column C)
# values for the rows whose figures were recalculated upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "C3"   = 0
    "C5"   = 0
    "C32"  = 6
    "C42"  = 5
    "C45"  = 0
    "C72"  = 6
    "C73"  = 28
    "C74"  = 0
    "C93"  = 1
    "C114" = 4
    "C125" = 0
    "C150" = 0
    "C188" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
